$wb = $excel.ActiveWorkbook

# --- Rename sheets to their new, more readable names --------------------
# (Excel auto-updates any sheet-qualified defined names, e.g. the hidden
# ExternalData_* ranges, when a sheet is renamed.)
$wsLiabs = $wb.Worksheets.Item("_liabs2")
$wsLiabs.Name = "_Liabilities"

$wsAssets = $wb.Worksheets.Item("_assets")
$wsAssets.Name = "_Assets"

$wsEquity = $wb.Worksheets.Item("_equity")
$wsEquity.Name = "_Equity"

$wsStructure = $wb.Worksheets.Item("_Structure")

# --- Shorten the "Short Name" labels (column F) on each data sheet ------

# _Liabilities
$wsLiabs.Range("F2").Value = "Insurance "
$wsLiabs.Range("F3").Value = "Investment - discretionary"
$wsLiabs.Range("F4").Value = "Investment - not  discretionary"
$wsLiabs.Range("F5").Value = "Unallocated surplus"
$wsLiabs.Range("F6").Value = "Third party interest"
$wsLiabs.Range("F7").Value = "Subordinated"
$wsLiabs.Range("F8").Value = "DB pension liability"
$wsLiabs.Range("F11").Value = "Derivative Liabilities"
$wsLiabs.Range("F12").Value = "Leases"
$wsLiabs.Range("F13").Value = "Other financial"
$wsLiabs.Range("F15").Value = "Accruals deferred income"
$wsLiabs.Range("F16").Value = "Liabilities for sale"

# _Assets
$wsAssets.Range("F2").Value = "Intangible"
$wsAssets.Range("F3").Value = "Deferred acquisitions"
$wsAssets.Range("F4").Value = "Joint ventures"
$wsAssets.Range("F5").Value = "Property plant"
$wsAssets.Range("F7").Value = "DB pension asset"
$wsAssets.Range("F9").Value = "Reinsurance"
$wsAssets.Range("F12").Value = "Equity securities"
$wsAssets.Range("F16").Value = "Accrued investment"
$wsAssets.Range("F17").Value = "Assets for sale"
$wsAssets.Range("F18").Value = "Cash"

# _Equity (text unchanged, kept for clarity / shared-string table reflow)
$wsEquity.Range("F4").Value = "Shares in employee trust"

# --- Restore each sheet's last-used selection ----------------------------
$wsLiabs.Activate()
$wsLiabs.Range("F20").Select()

$wsAssets.Activate()
$wsAssets.Range("F8").Select()

$wsEquity.Activate()
$wsEquity.Range("E2:E8").Select()

$wsStructure.Activate()
$wsStructure.Range("B2:B4").Select()
